$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 10.46510533333333
$ws.Range("H2").Value = 31.395316
$ws.Range("I2").Value = 0.5554075997074865
$ws.Range("J2").Value = 0.5554075997074865
$ws.Range("M2").Value = 0.034325
$ws.Range("N2").Value = 0.102975
$ws.Range("O2").Value = 0.004508979075184418
$ws.Range("P2").Value = 0.004508979075184418
$ws.Range("Q2").Value = 0.3592147405666667
$ws.Range("R2").Value = 3.2329326651
$ws.Range("S2").Value = 0.00250432124527946
$ws.Range("T2").Value = 0.00250432124527946
$ws.Range("G3").Value = 10.46510533333333
$ws.Range("H3").Value = 31.395316
$ws.Range("I3").Value = 0.5554075997074865
$ws.Range("J3").Value = 0.5554075997074865
$ws.Range("O3").Value = 0.9142039036746329
$ws.Range("P3").Value = 0.9142039036746329
$ws.Range("Q3").Value = 72.83145754454088
$ws.Range("R3").Value = 655.483117900868
$ws.Range("S3").Value = 0.5077557957831421
$ws.Range("T3").Value = 0.5077557957831421
$ws.Range("G4").Value = 10.46510533333333
$ws.Range("H4").Value = 31.395316
$ws.Range("I4").Value = 0.5554075997074865
$ws.Range("J4").Value = 0.5554075997074865
$ws.Range("M4").Value = 0.5818573333333333
$ws.Range("N4").Value = 1.745572
$ws.Range("O4").Value = 0.0764335772976724
$ws.Range("P4").Value = 0.0764335772976724
$ws.Range("Q4").Value = 6.089198282305778
$ws.Range("R4").Value = 54.80278454075201
$ws.Range("S4").Value = 0.04245178970395686
$ws.Range("T4").Value = 0.04245178970395686
$ws.Range("G5").Value = 10.46510533333333
$ws.Range("H5").Value = 31.395316
$ws.Range("I5").Value = 0.5554075997074865
$ws.Range("J5").Value = 0.5554075997074865
$ws.Range("M5").Value = 0.036948
$ws.Range("N5").Value = 0.110844
$ws.Range("O5").Value = 0.004853539952510238
$ws.Range("P5").Value = 0.004853539952510237
$ws.Range("Q5").Value = 0.386664711856
$ws.Range("R5").Value = 3.479982406704
$ws.Range("S5").Value = 0.002695692975108099
$ws.Range("T5").Value = 0.002695692975108099
$ws.Range("I6").Value = 0.3053945925621632
$ws.Range("J6").Value = 0.3053945925621632
$ws.Range("M6").Value = 0.034325
$ws.Range("N6").Value = 0.102975
$ws.Range("O6").Value = 0.004508979075184418
$ws.Range("P6").Value = 0.004508979075184418
$ws.Range("Q6").Value = 0.1975166335416667
$ws.Range("R6").Value = 1.777649701875
$ws.Range("S6").Value = 0.001377017827537265
$ws.Range("T6").Value = 0.001377017827537265
$ws.Range("I7").Value = 0.3053945925621632
$ws.Range("J7").Value = 0.3053945925621632
$ws.Range("O7").Value = 0.9142039036746329
$ws.Range("P7").Value = 0.9142039036746329
$ws.Range("S7").Value = 0.2791929286814536
$ws.Range("T7").Value = 0.2791929286814536
$ws.Range("I8").Value = 0.3053945925621632
$ws.Range("J8").Value = 0.3053945925621632
$ws.Range("M8").Value = 0.5818573333333333
$ws.Range("N8").Value = 1.745572
$ws.Range("O8").Value = 0.0764335772976724
$ws.Range("P8").Value = 0.0764335772976724
$ws.Range("Q8").Value = 3.348186502011111
$ws.Range("R8").Value = 30.13367851810001
$ws.Range("S8").Value = 0.02334240119689127
$ws.Range("T8").Value = 0.02334240119689127
$ws.Range("I9").Value = 0.3053945925621632
$ws.Range("J9").Value = 0.3053945925621632
$ws.Range("M9").Value = 0.036948
$ws.Range("N9").Value = 0.110844
$ws.Range("O9").Value = 0.004853539952510238
$ws.Range("P9").Value = 0.004853539952510237
$ws.Range("Q9").Value = 0.2126101843
$ws.Range("R9").Value = 1.9134916587
$ws.Range("S9").Value = 0.001482244856281045
$ws.Range("T9").Value = 0.001482244856281045
$ws.Range("G10").Value = 2.146766
$ws.Range("H10").Value = 6.440298
$ws.Range("I10").Value = 0.1139338891693565
$ws.Range("J10").Value = 0.1139338891693565
$ws.Range("M10").Value = 0.034325
$ws.Range("N10").Value = 0.102975
$ws.Range("O10").Value = 0.004508979075184418
$ws.Range("P10").Value = 0.004508979075184418
$ws.Range("Q10").Value = 0.07368774295
$ws.Range("R10").Value = 0.66318968655
$ws.Range("S10").Value = 0.0005137255222190092
$ws.Range("T10").Value = 0.0005137255222190093
$ws.Range("G11").Value = 2.146766
$ws.Range("H11").Value = 6.440298
$ws.Range("I11").Value = 0.1139338891693565
$ws.Range("J11").Value = 0.1139338891693565
$ws.Range("O11").Value = 0.9142039036746329
$ws.Range("P11").Value = 0.9142039036746329
$ws.Range("Q11").Value = 14.94032709723933
$ws.Range("R11").Value = 134.462943875154
$ws.Range("S11").Value = 0.1041588062394587
$ws.Range("T11").Value = 0.1041588062394587
$ws.Range("G12").Value = 2.146766
$ws.Range("H12").Value = 6.440298
$ws.Range("I12").Value = 0.1139338891693565
$ws.Range("J12").Value = 0.1139338891693565
$ws.Range("M12").Value = 0.5818573333333333
$ws.Range("N12").Value = 1.745572
$ws.Range("O12").Value = 0.0764335772976724
$ws.Range("P12").Value = 0.0764335772976724
$ws.Range("Q12").Value = 1.249111540050667
$ws.Range("R12").Value = 11.242003860456
$ws.Range("S12").Value = 0.008708374724650452
$ws.Range("T12").Value = 0.008708374724650454
$ws.Range("G13").Value = 2.146766
$ws.Range("H13").Value = 6.440298
$ws.Range("I13").Value = 0.1139338891693565
$ws.Range("J13").Value = 0.1139338891693565
$ws.Range("M13").Value = 0.036948
$ws.Range("N13").Value = 0.110844
$ws.Range("O13").Value = 0.004853539952510238
$ws.Range("P13").Value = 0.004853539952510237
$ws.Range("Q13").Value = 0.079318710168
$ws.Range("R13").Value = 0.713868391512
$ws.Range("S13").Value = 0.0005529826830283454
$ws.Range("T13").Value = 0.0005529826830283453
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.4760280000000001
$ws.Range("H14").Value = 1.428084
$ws.Range("I14").Value = 0.02526391856099382
$ws.Range("J14").Value = 0.02526391856099382
$ws.Range("M14").Value = 0.034325
$ws.Range("N14").Value = 0.102975
$ws.Range("O14").Value = 0.004508979075184418
$ws.Range("P14").Value = 0.004508979075184418
$ws.Range("Q14").Value = 0.0163396611
$ws.Range("R14").Value = 0.1470569499
$ws.Range("S14").Value = 0.0001139144801486844
$ws.Range("T14").Value = 0.0001139144801486844
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.4760280000000001
$ws.Range("H15").Value = 1.428084
$ws.Range("I15").Value = 0.02526391856099382
$ws.Range("J15").Value = 0.02526391856099382
$ws.Range("O15").Value = 0.9142039036746329
$ws.Range("P15").Value = 0.9142039036746329
$ws.Range("Q15").Value = 3.312896714148
$ws.Range("R15").Value = 29.816070427332
$ws.Range("S15").Value = 0.02309637297057857
$ws.Range("T15").Value = 0.02309637297057857
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.4760280000000001
$ws.Range("H16").Value = 1.428084
$ws.Range("I16").Value = 0.02526391856099382
$ws.Range("J16").Value = 0.02526391856099382
$ws.Range("M16").Value = 0.5818573333333333
$ws.Range("N16").Value = 1.745572
$ws.Range("O16").Value = 0.0764335772976724
$ws.Range("P16").Value = 0.0764335772976724
$ws.Range("Q16").Value = 0.276980382672
$ws.Range("R16").Value = 2.492823444048001
$ws.Range("S16").Value = 0.001931011672173822
$ws.Range("T16").Value = 0.001931011672173822
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.4760280000000001
$ws.Range("H17").Value = 1.428084
$ws.Range("I17").Value = 0.02526391856099382
$ws.Range("J17").Value = 0.02526391856099382
$ws.Range("M17").Value = 0.036948
$ws.Range("N17").Value = 0.110844
$ws.Range("O17").Value = 0.004853539952510238
$ws.Range("P17").Value = 0.004853539952510237
$ws.Range("Q17").Value = 0.017588282544
$ws.Range("R17").Value = 0.158294542896
$ws.Range("S17").Value = 0.0001226194380927485
$ws.Range("T17").Value = 0.0001226194380927484
